# Auto-generated edit script: updates cryptos list values per diff
# Applies updated prices, volume-change percentages, and (for rows 43/44)
# swaps the ARBITRUM / RenderToken entries to match the new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps storing values as plain text, exactly as
# the source data feed provides them (prevents Excel from re-interpreting
# strings such as '212.54' or '0.998' as numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '29.761.37'
$ws.Cells.Item(2, 5).Value = '  -0.24%  '
$ws.Cells.Item(3, 4).Value = '1.603.36'
$ws.Cells.Item(3, 5).Value = '  -0.99%  '
$ws.Cells.Item(4, 4).Value = '0.998'
$ws.Cells.Item(4, 5).Value = '  +0.35%  '
$ws.Cells.Item(5, 4).Value = '212.54'
$ws.Cells.Item(5, 5).Value = '  -1.10%  '
$ws.Cells.Item(6, 5).Value = '  +0.06%  '
$ws.Cells.Item(7, 5).Value = '  +0.35%  '
$ws.Cells.Item(8, 4).Value = '28.20'
$ws.Cells.Item(8, 5).Value = '  +4.15%  '
$ws.Cells.Item(9, 5).Value = '  +0.94%  '
$ws.Cells.Item(10, 4).Value = '0.0604'
$ws.Cells.Item(10, 5).Value = '  +0.50%  '
$ws.Cells.Item(11, 4).Value = '0.0909'
$ws.Cells.Item(11, 5).Value = '  -0.73%  '
$ws.Cells.Item(12, 4).Value = '1.832.72'
$ws.Cells.Item(12, 5).Value = '  -0.72%  '
$ws.Cells.Item(13, 4).Value = '1.603.13'
$ws.Cells.Item(13, 5).Value = '  -0.98%  '
$ws.Cells.Item(14, 5).Value = '  +1.10%  '
$ws.Cells.Item(15, 4).Value = '29.759.49'
$ws.Cells.Item(15, 5).Value = '  -0.29%  '
$ws.Cells.Item(16, 5).Value = '  -0.40%  '
$ws.Cells.Item(17, 4).Value = '64.17'
$ws.Cells.Item(17, 5).Value = '  +0.57%  '
$ws.Cells.Item(18, 4).Value = '241.85'
$ws.Cells.Item(18, 5).Value = '  -1.92%  '
$ws.Cells.Item(19, 5).Value = '  +2.73%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0698'
$ws.Cells.Item(20, 5).Value = '  -0.13%  '
$ws.Cells.Item(21, 5).Value = '  +0.37%  '
$ws.Cells.Item(22, 5).Value = '  -1.35%  '
$ws.Cells.Item(23, 4).Value = '9.41'
$ws.Cells.Item(23, 5).Value = '  +1.02%  '
$ws.Cells.Item(24, 4).Value = '2.12'
$ws.Cells.Item(24, 5).Value = '  -0.48%  '
$ws.Cells.Item(25, 4).Value = '155.20'
$ws.Cells.Item(25, 5).Value = '  -0.65%  '
$ws.Cells.Item(26, 4).Value = '15.44'
$ws.Cells.Item(26, 5).Value = '  -0.05%  '
$ws.Cells.Item(27, 4).Value = '0.110'
$ws.Cells.Item(27, 5).Value = '  +0.69%  '
$ws.Cells.Item(28, 5).Value = '  +0.19%  '
$ws.Cells.Item(29, 5).Value = '  +0.35%  '
$ws.Cells.Item(30, 4).Value = '0.0478'
$ws.Cells.Item(30, 5).Value = '  +0.55%  '
$ws.Cells.Item(31, 5).Value = '  -0.09%  '
$ws.Cells.Item(32, 5).Value = '  -0.76%  '
$ws.Cells.Item(33, 5).Value = '  +1.93%  '
$ws.Cells.Item(34, 4).Value = '1.421.31'
$ws.Cells.Item(34, 5).Value = '  -1.97%  '
$ws.Cells.Item(35, 5).Value = '  +2.34%  '
$ws.Cells.Item(36, 5).Value = '  +1.66%  '
$ws.Cells.Item(37, 5).Value = '  -1.94%  '
$ws.Cells.Item(38, 5).Value = '  -0.49%  '
$ws.Cells.Item(39, 5).Value = '  +0.77%  '
$ws.Cells.Item(40, 4).Value = '0.545'
$ws.Cells.Item(40, 5).Value = '  +0.50%  '
$ws.Cells.Item(41, 4).Value = '56.31'
$ws.Cells.Item(41, 5).Value = '  -1.42%  '
$ws.Cells.Item(42, 4).Value = '0.0495'
$ws.Cells.Item(42, 5).Value = '  +5.73%  '
$ws.Cells.Item(43, 2).Value = 'ARBITRUM'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(43, 4).Value = '0.817'
$ws.Cells.Item(43, 5).Value = '  +1.48%  '
$ws.Cells.Item(44, 2).Value = 'RenderToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(44, 4).Value = '1.95'
$ws.Cells.Item(44, 5).Value = '  -1.04%  '
$ws.Cells.Item(45, 5).Value = '  +0.34%  '
$ws.Cells.Item(46, 4).Value = '67.33'
$ws.Cells.Item(46, 5).Value = '  -3.27%  '
$ws.Cells.Item(47, 4).Value = '0.984'
$ws.Cells.Item(47, 5).Value = '  +17.54%  '
$ws.Cells.Item(48, 5).Value = '  +0.45%  '
$ws.Cells.Item(49, 4).Value = '1.741.66'
$ws.Cells.Item(49, 5).Value = '  -1.02%  '
$ws.Cells.Item(50, 4).Value = '86.55'
$ws.Cells.Item(50, 5).Value = '  -0.21%  '
$ws.Cells.Item(51, 4).Value = '0.0₆0103'
$ws.Cells.Item(51, 5).Value = '  +2.01%  '
